# Update "想去人数" (F column) values across the sheets, and consolidate the
# duplicate "第二届北极光动漫展" row on "全部类型" with the "万圣漫控嘉年华10"
# row that follows it (row 19 merges into row 18, row 19 is removed).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 610
$ws1.Range("F3").Value = 205
$ws1.Range("F4").Value = 514
$ws1.Range("F5").Value = 510
$ws1.Range("F6").Value = 289
$ws1.Range("F7").Value = 2603
$ws1.Range("F8").Value = 443
$ws1.Range("F9").Value = 7108
$ws1.Range("F10").Value = 188
$ws1.Range("F11").Value = 447
$ws1.Range("F12").Value = 11
$ws1.Range("F13").Value = 120
$ws1.Range("F14").Value = 37

# --- Sheet 2: 演出 -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 13
$ws2.Range("F3").Value = 17
$ws2.Range("F4").Value = 1
$ws2.Range("F5").Value = 1

# --- Sheet 4: 全部类型 --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 610
$ws4.Range("F3").Value = 205
$ws4.Range("F4").Value = 514
$ws4.Range("F5").Value = 510
$ws4.Range("F6").Value = 289
$ws4.Range("F7").Value = 13
$ws4.Range("F8").Value = 17
$ws4.Range("F9").Value = 2603
$ws4.Range("F10").Value = 443
$ws4.Range("F11").Value = 7108
$ws4.Range("F12").Value = 188
$ws4.Range("F13").Value = 447
$ws4.Range("F14").Value = 11
$ws4.Range("F15").Value = 1
$ws4.Range("F16").Value = 1
$ws4.Range("F17").Value = 120

# Row 18 currently duplicates row 17 (南宁·第二届北极光动漫展). Overwrite it
# with row 19's event (南宁·万圣漫控嘉年华10), then delete row 19 so the
# table ends at row 18.
# Force the date-like text to stay a plain string (not auto-converted to a
# date serial) by switching the cell to text format before assigning, then
# drop back to the default style so no stray formatting is left behind.
$ws4.Range("B18").NumberFormat = "@"
$ws4.Range("B18").Value = "2024-11-02"
$ws4.Range("B18").Style = "Normal"
$ws4.Range("C18").Value = "南宁·万圣漫控嘉年华10"
$ws4.Range("D18").Value = "亭洪路45号 百益上河城"
$ws4.Range("E18").Value = "2024.11.02 11:00-11.03 22:00"
$ws4.Range("F18").Value = 37
$ws4.Range("G18").Value = 50
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=87820"
$ws4.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg"

$ws4.Rows.Item(19).Delete()
